$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.929.77"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.464.54"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "2.464.80"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.112"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "2.911.24"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "62.831.50"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "2.466.76"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +19.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "656.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Value = "0.0₃0985"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").Value = "2.584.85"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("E30").Value = "  -14.40%  "
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.135"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.90%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "151.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("D44").Value = "0.0₆0310"
$ws.Range("E44").Value = "  -41.12%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "153.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.608"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0512"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.80%  "
